$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values that look numeric (e.g. "1.00", "575.74") must be forced to
# plain text so Excel does not silently coerce them into Number cells and
# normalize/round their display (the source data keeps them as literal text,
# same as values like "62.014.73" which are not valid numbers anyway).
$textCells = @("D4", "D5", "D6", "D8", "D9", "D11", "D13", "D14", "D18", "D19", "D20", "D21", "D22", "D23", "D26", "D28", "D30", "D31", "D33", "D34", "D35", "D36", "D38", "D39", "D41", "D42", "D43", "D44", "D46", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.014.73"
$ws.Range("E2").Value = "  -1.44%  "

$ws.Range("D3").Value = "3.411.16"
$ws.Range("E3").Value = "  -1.75%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "575.74"
$ws.Range("E5").Value = "  -0.33%  "

$ws.Range("D6").Value = "148.21"
$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "0.485"
$ws.Range("E8").Value = "  +1.08%  "

$ws.Range("D9").Value = "7.98"
$ws.Range("E9").Value = "  +4.31%  "

$ws.Range("E10").Value = "  -1.02%  "

$ws.Range("D11").Value = "0.414"
$ws.Range("E11").Value = "  +2.70%  "

$ws.Range("D12").Value = "3.996.73"
$ws.Range("E12").Value = "  -1.80%  "

$ws.Range("D13").Value = "0.129"
$ws.Range("E13").Value = "  +0.15%  "

$ws.Range("D14").Value = "28.32"
$ws.Range("E14").Value = "  -5.20%  "

$ws.Range("D15").Value = "3.396.67"
$ws.Range("E15").Value = "  -2.35%  "

$ws.Range("E16").Value = "  -0.11%  "

$ws.Range("D17").Value = "62.040.49"
$ws.Range("E17").Value = "  -1.50%  "

$ws.Range("D18").Value = "6.40"
$ws.Range("E18").Value = "  +0.90%  "

$ws.Range("D19").Value = "14.46"
$ws.Range("E19").Value = "  +0.48%  "

$ws.Range("D20").Value = "8.95"
$ws.Range("E20").Value = "  -2.89%  "

$ws.Range("D21").Value = "380.06"
$ws.Range("E21").Value = "  -2.42%  "

$ws.Range("D22").Value = "0.565"
$ws.Range("E22").Value = "  +1.17%  "

$ws.Range("D23").Value = "74.78"
$ws.Range("E23").Value = "  +0.14%  "

$ws.Range("E24").Value = "  -0.07%  "

$ws.Range("D25").Value = "3.573.74"
$ws.Range("E25").Value = "  -0.90%  "

$ws.Range("D26").Value = "0.0000112"
$ws.Range("E26").Value = "  -2.51%  "

$ws.Range("E27").Value = "  +1.02%  "

$ws.Range("D28").Value = "7.64"
$ws.Range("E28").Value = "  +0.99%  "

$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "2.13"
$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "7.90"
$ws.Range("E31").Value = "  -3.00%  "

$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("D33").Value = "1.34"
$ws.Range("E33").Value = "  -1.75%  "

$ws.Range("D34").Value = "23.09"
$ws.Range("E34").Value = "  -2.44%  "

$ws.Range("D35").Value = "5.46"
$ws.Range("E35").Value = "  +3.78%  "

$ws.Range("D36").Value = "1.62"
$ws.Range("E36").Value = "  +3.69%  "

$ws.Range("E37").Value = "  -2.15%  "

$ws.Range("D38").Value = "169.22"
$ws.Range("E38").Value = "  -0.74%  "

$ws.Range("D39").Value = "30.45"
$ws.Range("E39").Value = "  -4.82%  "

$ws.Range("D40").Value = "3.446.24"
$ws.Range("E40").Value = "  -1.90%  "

$ws.Range("D41").Value = "0.0781"
$ws.Range("E41").Value = "  +3.46%  "

$ws.Range("D42").Value = "0.784"
$ws.Range("E42").Value = "  -1.75%  "

$ws.Range("D43").Value = "42.42"
$ws.Range("E43").Value = "  +0.06%  "

$ws.Range("D44").Value = "4.36"
$ws.Range("E44").Value = "  -2.36%  "

$ws.Range("E45").Value = "  -2.33%  "

$ws.Range("D46").Value = "1.17"
$ws.Range("E46").Value = "  -2.92%  "

$ws.Range("D47").Value = "2.541.94"
$ws.Range("E47").Value = "  -2.69%  "

$ws.Range("D48").Value = "6.89"
$ws.Range("E48").Value = "  +2.57%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "22.68"
$ws.Range("E49").Value = "  -1.09%  "

$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "2.20"
$ws.Range("E50").Value = "  -3.17%  "

$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  -0.08%  "

# Restore the default "Normal" cell style on the price cells now that the
# text value is committed, so no stray number-format style lingers on them.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
